$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

# Insert a blank row below the "vac" linked_table query (row 2), pushing the
# "linked_visitdate" query down to row 4 and leaving a blank separator row 3.
$ws.Rows.Item(3).Insert()

# Row 2 ("vac" query): the selection clause and the newRow initialisation map
# now key off REGIDC instead of the (non-existent) generic ID field, and the
# initial map's visit-id key is renamed VISITIDC.
$ws.Range("E2").Value = "REGIDC = ?"
$ws.Range("G2").Value = "{REGIDC: data('REGIDC'), VISITDATE: data('CONT'), VISITIDC: data('_id')}"

# Row 4 ("linked_visitdate" query, previously row 3): fix the broken link -
# it used to select on a non-existent "ID" field; it should join through
# REGIDC like the other linked-table query does.
$ws.Range("E4").Value = "REGIDC = ?"
$ws.Range("F4").Value = "[data('REGIDC')]"

# Column widths shift slightly now that the content (and therefore Excel's
# best-fit measurement) has changed.
$ws.Columns.Item(1).ColumnWidth = 14.6
$ws.Columns.Item(5).ColumnWidth = 9.45
$ws.Columns.Item(6).ColumnWidth = 14.25

$ws.Range("G2").Select()
